$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Minimize the workbook window (workbookView minimized="1" in the target XML)
$wb.Windows.Item(1).WindowState = -4140

# "Last job" / "Last location" values updated for the ejobs row
$ws.Range("C2").Value = "C++"
$ws.Range("D2").Value = "Cluj"

# New "Email" column: copy the existing header formatting onto E1, then set its text
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("E1").Value = "Email"

# New contact e-mail address, wired up as a mailto hyperlink (gets the workbook's
# Hyperlink style automatically, same as the other link cells in column A)
$ws.Range("E2").Value = "catalinnm99@gmail.com"
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:catalinnm99@gmail.com", "", "", "catalinnm99@gmail.com")

# Selection ends up on E3 before the file is saved
$ws.Range("E3").Select()
